# "more phone number validation"
# Apply a Text number format to the phoneNumber column (F2:F12) and store
# the phone numbers as validated text strings (preserving leading zeros
# and "+" country-code prefixes) instead of raw numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text format ("@") to the whole phone-number data range first so
# every value below is (re)written as text.
$ws.Range("F2:F12").NumberFormat = "@"

$ws.Range("F2").Value = "0765893872"
$ws.Range("F3").Value = "0721748392"
$ws.Range("F4").Value = "0721745592"
$ws.Range("F5").Value = "0712348392"
$ws.Range("F6").Value = "0721799423"
$ws.Range("F7").Value = "+40765126291"
$ws.Range("F8").Value = "wd"
$ws.Range("F9").Value = "wd"
$ws.Range("F10").Value = "a"
$ws.Range("F11").Value = "+765472891"
$ws.Range("F12").Value = "awd"

$ws.Range("F13").Select() | Out-Null
